$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure dollar-amount strings are stored as plain text (not auto-converted
# to currency numbers by Excel's input parsing) - matches existing cells like F10:Q10.
$ws.Range("F4:M4").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"

# New Glass Stop column (M) values
$ws.Range("M1").Value = "Glass Stop (Use with BE9-2514, BE9-2515, BE9-2517)"
$ws.Range("M2").Value = "E9-2519"
$ws.Range("M3").Value = "24.0 ft"
$ws.Range("M4").Value = "$85.05"

# Updated prices in row 4 (F4:L4) reflecting the recalculated costs
$ws.Range("F4").Value = "$279.33"
$ws.Range("G4").Value = "$209.50"
$ws.Range("H4").Value = "$150.00"
$ws.Range("I4").Value = "$497.00"
$ws.Range("J4").Value = "$263.50"
$ws.Range("K4").Value = "$222.00"
$ws.Range("L4").Value = "$132.00"

# Updated Grand Total
$ws.Range("E14").Value = "$2154.48"
